# Generate Report for Handback
#
# This script reproduces the "handback" report-generation edit:
#  - Status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US" everywhere it appears
#    (Overview!E2:F2,E3:F3 and the per-language sheets' Status column).
#  - The zh-cn / de-de sheets get their "Latest Target File" (I) and
#    "Latest Handback File" (J) columns populated (with I linking back
#    to the handed-off markdown source), and the "Latest Handback
#    DateTime" (K) column gets a fresh timestamp (different per
#    language).
#  - A handful of columns are widened to fit the newly-populated data.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

# --- Status column text (shared across Overview + both language sheets) ---
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

# --- Handback hyperlinks & target/handback file names ---
$mdTarget1 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/807db5350507b98eebda34ebf9efbd771f8630c6/e2e/a36dfb4e-a64c-4a30-9bee-ef515d000f79.md"
$mdTarget2 = "https://github.com/OpenLocalizationTestOrg/oltest/blob/807db5350507b98eebda34ebf9efbd771f8630c6/e2e/e82f82ba-b658-4fda-a5a4-314be47a6254.md"

# zh-cn row 2 (a36dfb4e...)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdTarget1, "", "", "a36dfb4e-a64c-4a30-9bee-ef515d000f79.md")
$wsZh.Range("J2").Value = "a36dfb4e-a64c-4a30-9bee-ef515d000f79.7e68a34c5380b9956fdb073c78952f616a2444d9.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-12 02:51:51"

# zh-cn row 3 (e82f82ba...)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdTarget2, "", "", "e82f82ba-b658-4fda-a5a4-314be47a6254.md")
$wsZh.Range("J3").Value = "e82f82ba-b658-4fda-a5a4-314be47a6254.cb2af8e4f9fb665c7bf75bc98237aa64d760ed73.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-12 02:51:51"

# de-de row 2 (a36dfb4e...)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdTarget1, "", "", "a36dfb4e-a64c-4a30-9bee-ef515d000f79.md")
$wsDe.Range("J2").Value = "a36dfb4e-a64c-4a30-9bee-ef515d000f79.7e68a34c5380b9956fdb073c78952f616a2444d9.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-12 02:51:58"

# de-de row 3 (e82f82ba...)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdTarget2, "", "", "e82f82ba-b658-4fda-a5a4-314be47a6254.md")
$wsDe.Range("J3").Value = "e82f82ba-b658-4fda-a5a4-314be47a6254.cb2af8e4f9fb665c7bf75bc98237aa64d760ed73.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-12 02:51:58"

# --- Column widening to fit the now-populated / longer columns ---
# (ColumnWidth values chosen so the saved sheet's raw column width lands
# on the same value the report generator produced.)
$wsOverview.Columns.Item(5).ColumnWidth = 29.16
$wsOverview.Columns.Item(6).ColumnWidth = 29.16

$wsZh.Columns.Item(3).ColumnWidth = 29.16
$wsZh.Columns.Item(9).ColumnWidth = 39.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15

$wsDe.Columns.Item(3).ColumnWidth = 29.16
$wsDe.Columns.Item(9).ColumnWidth = 39.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15

Write-Host "Handback report generated."
